# Update "想去人数" (interested-in count) figures on the sheets that list
# the full con data: "展览" and "全部类型" both received the same updates.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 10
    $ws.Range("F4").Value = 953
    $ws.Range("F6").Value = 433
}
